$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the measured values (row 4-8, columns C/D/E)
$ws.Range("C4").Value = 403
$ws.Range("D4").Value = 1073
$ws.Range("E4").Value = 1.3

$ws.Range("C5").Value = 401
$ws.Range("D5").Value = 671
$ws.Range("E5").Value = 2

$ws.Range("C6").Value = 208
$ws.Range("D6").Value = 671
$ws.Range("E6").Value = 1

$ws.Range("C7").Value = 111
$ws.Range("D7").Value = 469
$ws.Range("E7").Value = 0.8

$ws.Range("C8").Value = 110
$ws.Range("D8").Value = 301
$ws.Range("E8").Value = 1.2

# Add new label "h" in G16
$ws.Range("G16").Value = "h"

# Update selection / view (topLeftCell moved to B1, active cell to B10)
$excel.ActiveWindow.DisplayGridlines = $true
$ws.Range("B10").Select()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
